$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 11, shifting existing rows 11-22 down to 12-23.
$ws.Rows(11).Insert()

# Populate the newly inserted row 11 with the new weekly price entry.
$ws.Range("A11").Value = 11
$ws.Range("B11").Value = "Vega Monumental Concepción"
$ws.Range("C11").Value = "Bíobío"
$ws.Range("D11").Value = 45274
$ws.Range("E11").Value = 8
$ws.Range("F11").Value = "Fruta"
$ws.Range("G11").Value = 100101
$ws.Range("H11").Value = "Berries"
$ws.Range("I11").Value = 100101004
$ws.Range("J11").Value = "Frambuesa"
$ws.Range("K11").Value = "Sin especificar"
$ws.Range("L11").Value = "Primera"
$ws.Range("M11").Value = 100
$ws.Range("N11").Value = 8000
$ws.Range("O11").Value = 9000
$ws.Range("P11").Value = 8500
$ws.Range("Q11").Value = "$/bandeja 2 kilos"
$ws.Range("R11").Value = "Región de Ñuble"
$ws.Range("S11").Value = 4250
$ws.Range("T11").Value = 2
